$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.133.26"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").Value = "2.367.23"
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.65"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.20"
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").Value = "2.372.05"
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.83"
$ws.Range("E12").Value = "  +4.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.322"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "2.787.91"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "56.117.11"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.42"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "2.293.14"
$ws.Range("E18").Value = "  -7.82%  "
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.03"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "306.91"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.63"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("E27").Value = "  -5.73%  "
$ws.Range("E28").Value = "  -4.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.39"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").Value = "0.0₃0710"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -6.80%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.61"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("E37").Value = "  -6.19%  "
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.08"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("E41").Value = "  -5.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.07"
$ws.Range("E42").Value = "  -5.83%  "
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.70"
$ws.Range("E44").Value = "  -4.33%  "
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0902"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "239.10"
$ws.Range("E47").Value = "  -6.90%  "
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.04"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  -0.77%  "
